# Update the workbook to add 2021 and 2022 data rows and correct the
# 2019 manufacturing-investment growth figure (B3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing 2019 value (B3): 3.07744464978055 -> 3.1 ---
$ws.Range("B3").Value = 3.1

# --- Make new label cells (A5, A6) match the formatting already used
#     for the other year labels in column A (A2:A4), e.g. bold/centered
#     style, by copying the format from A4 before writing into them. ---
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 5: 2021年 ---
$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 13.5
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 4.4

# --- Row 6: 2022年 (real-estate figure not yet available -> left blank) ---
$ws.Range("A6").Value = "2022年"
$ws.Range("B6").Value = 9.1
$ws.Range("C6").Value = 9.4
$ws.Range("D6").Value = ""
